$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel stores them as literal text (matching original inlineStr formatting)
# instead of converting them to floating point numbers.
$textForceCells = @("D5", "D10", "D11", "D16", "D20", "D21", "D22", "D25", "D27", "D28", "D29", "D33", "D39", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply all cell value updates described by the diff
$ws.Range("D2").Value = "26.637.91"
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("D3").Value = "1.588.18"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("D5").Value = "210.76"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "19.55"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").Value = "0.0833"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "1.810.72"
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("D13").Value = "1.583.63"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("E15").Value = "  -4.40%  "
$ws.Range("D16").Value = "64.74"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "26.620.00"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D20").Value = "208.06"
$ws.Range("E20").Value = "  -4.21%  "
$ws.Range("D21").Value = "6.71"
$ws.Range("E21").Value = "  -3.32%  "
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("E24").Value = "  -2.24%  "
$ws.Range("D25").Value = "146.67"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").Value = "7.22"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").Value = "0.113"
$ws.Range("E28").Value = "  -3.88%  "
$ws.Range("D29").Value = "15.27"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("E32").Value = "  -4.26%  "
$ws.Range("D33").Value = "0.660"
$ws.Range("E33").Value = "  +19.80%  "
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("D35").Value = "1.305.74"
$ws.Range("E35").Value = "  -3.33%  "
$ws.Range("E37").Value = "  -5.36%  "
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").Value = "0.827"
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("D44").Value = "62.57"
$ws.Range("E44").Value = "  -4.64%  "
$ws.Range("D45").Value = "1.724.55"
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("D46").Value = "89.67"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "0.837"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0503"
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0978"
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "7.51"
$ws.Range("E51").Value = "  -1.25%  "
